$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '20.224.05'
$ws.Range('E2').Value = '  +2.18%  '
$ws.Range('D3').Value = '1.442.70'
$ws.Range('E3').Value = '  +4.14%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.9118'
$ws.Range('E5').Value = '  -9.23%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '277.95'
$ws.Range('E6').Value = '  +3.68%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3650'
$ws.Range('E7').Value = '  +1.00%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3111'
$ws.Range('E8').Value = '  +3.81%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '38.98'
$ws.Range('E9').Value = '  +0.48%  '
$ws.Range('E10').Value = '  +7.07%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.06523'
$ws.Range('E11').Value = '  +3.36%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  -0.60%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.390'
$ws.Range('E13').Value = '  +3.31%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '17.61'
$ws.Range('E14').Value = '  +8.77%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.051'
$ws.Range('E15').Value = '  +1.26%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.00001017'
$ws.Range('E16').Value = '  +4.39%  '
$ws.Range('D17').Value = '1.442.41'
$ws.Range('E17').Value = '  +3.84%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.9396'
$ws.Range('E18').Value = '  -6.45%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.05635'
$ws.Range('E19').Value = '  +0.48%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '68.06'
$ws.Range('E20').Value = '  -2.12%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.383'
$ws.Range('E22').Value = '  -1.29%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.80'
$ws.Range('E23').Value = '  +3.67%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.262'
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('D25').Value = '20.262.62'
$ws.Range('E25').Value = '  +2.27%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.164'
$ws.Range('E26').Value = '  +2.70%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '137.95'
$ws.Range('E27').Value = '  +2.43%  '
$ws.Range('E28').Value = '  +3.48%  '
$ws.Range('D29').Value = '1.595.72'
$ws.Range('E29').Value = '  +2.70%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '109.99'
$ws.Range('E30').Value = '  +2.62%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.827'
$ws.Range('E31').Value = '  +0.24%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.8046'
$ws.Range('E32').Value = '  +2.87%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.826'
$ws.Range('E33').Value = '  -6.63%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.07699'
$ws.Range('E34').Value = '  +2.11%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.05919'
$ws.Range('E35').Value = '  +7.66%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.446'
$ws.Range('E36').Value = '  +11.79%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.670'
$ws.Range('E37').Value = '  +1.04%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.136'
$ws.Range('E38').Value = '  +10.67%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01991'
$ws.Range('E39').Value = '  +0.78%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '10.14'
$ws.Range('E40').Value = '  +2.47%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.9331'
$ws.Range('E41').Value = '  -7.10%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.1840'
$ws.Range('E42').Value = '  -1.30%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '7.157'
$ws.Range('E43').Value = '  -12.64%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.521'
$ws.Range('E44').Value = '  +2.04%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.5228'
$ws.Range('E45').Value = '  +2.28%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '12.06'
$ws.Range('E46').Value = '  +3.68%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '119.11'
$ws.Range('E47').Value = '  +11.21%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.5133'
$ws.Range('E48').Value = '  +5.09%  '
$ws.Range('E49').Value = '  +3.71%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06350'
$ws.Range('E50').Value = '  +5.00%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.9931'
$ws.Range('E51').Value = '  -1.09%  '
